# Scheduled runner update: refresh computed market-price / profit columns
# (H:N) on several rows across the per-job sheets (ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, WVR), matching newly recalculated Universalis price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 140742110
$ws.Cells.Item(62, 9).Value = 158334370
$ws.Cells.Item(62, 11).Value = 158334370
$ws.Cells.Item(62, 13).Value = -158333746
$ws.Cells.Item(65, 8).Value = 140742110
$ws.Cells.Item(65, 9).Value = 158334370
$ws.Cells.Item(65, 11).Value = 791671850
$ws.Cells.Item(65, 13).Value = -791668730
$ws.Cells.Item(80, 8).Value = 1239.9667
$ws.Cells.Item(80, 9).Value = 1141.0588
$ws.Cells.Item(80, 10).Value = 1369.3077
$ws.Cells.Item(80, 11).Value = 3423.1764
$ws.Cells.Item(80, 12).Value = 4107.9231
$ws.Cells.Item(80, 13).Value = -2425.1764
$ws.Cells.Item(80, 14).Value = -6103.9231
$ws.Cells.Item(83, 8).Value = 1239.9667
$ws.Cells.Item(83, 9).Value = 1141.0588
$ws.Cells.Item(83, 10).Value = 1369.3077
$ws.Cells.Item(83, 11).Value = 10269.5292
$ws.Cells.Item(83, 12).Value = 12323.7693
$ws.Cells.Item(83, 13).Value = -5277.529200000001
$ws.Cells.Item(83, 14).Value = -22307.7693
$ws.Cells.Item(141, 8).Value = 999.625
$ws.Cells.Item(141, 9).Value = 999.625
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 2998.875
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 13).Value = 2181.125
$ws.Cells.Item(141, 14).Value = $null
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 8206.571
$ws.Cells.Item(45, 9).Value = 10852.9
$ws.Cells.Item(45, 10).Value = 1590.75
$ws.Cells.Item(45, 11).Value = 10852.9
$ws.Cells.Item(45, 12).Value = 1590.75
$ws.Cells.Item(45, 13).Value = -10475.9
$ws.Cells.Item(45, 14).Value = -2344.75
$ws.Cells.Item(50, 8).Value = 1126.1666
$ws.Cells.Item(50, 9).Value = 285
$ws.Cells.Item(50, 10).Value = 1294.4
$ws.Cells.Item(50, 11).Value = 285
$ws.Cells.Item(50, 12).Value = 1294.4
$ws.Cells.Item(50, 13).Value = 429
$ws.Cells.Item(50, 14).Value = -2722.4
$ws.Cells.Item(132, 8).Value = 5474.9473
$ws.Cells.Item(132, 9).Value = 4207.4116
$ws.Cells.Item(132, 11).Value = 12622.2348
$ws.Cells.Item(132, 13).Value = -10092.2348
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(10, 8).Value = 1343.1428
$ws.Cells.Item(10, 10).Value = 1343.1428
$ws.Cells.Item(10, 12).Value = 1343.1428
$ws.Cells.Item(10, 14).Value = -1623.1428
$ws.Cells.Item(82, 8).Value = 27996.8
$ws.Cells.Item(82, 9).Value = 13842.538
$ws.Cells.Item(82, 11).Value = 13842.538
$ws.Cells.Item(82, 13).Value = -13459.538
$ws.Cells.Item(85, 8).Value = 27996.8
$ws.Cells.Item(85, 9).Value = 13842.538
$ws.Cells.Item(85, 11).Value = 13842.538
$ws.Cells.Item(85, 13).Value = -12516.538
$ws.Cells.Item(107, 8).Value = 1968.0555
$ws.Cells.Item(107, 9).Value = 1906.5883
$ws.Cells.Item(107, 10).Value = 3013
$ws.Cells.Item(107, 11).Value = 1906.5883
$ws.Cells.Item(107, 12).Value = 3013
$ws.Cells.Item(107, 13).Value = 13.41170000000011
$ws.Cells.Item(107, 14).Value = -6853
$ws.Cells.Item(109, 8).Value = 99993
$ws.Cells.Item(109, 10).Value = 99993
$ws.Cells.Item(109, 12).Value = 99993
$ws.Cells.Item(109, 14).Value = -102767
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(60, 8).Value = 60000
$ws.Cells.Item(60, 9).Value = 0
$ws.Cells.Item(60, 11).Value = 0
$ws.Cells.Item(60, 13).Value = $null
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).Value = $null
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).Value = $null
$ws.Cells.Item(134, 8).Value = 11392.5
$ws.Cells.Item(134, 9).Value = 11105.454
$ws.Cells.Item(134, 11).Value = 33316.362
$ws.Cells.Item(134, 13).Value = -30781.362
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 3611
$ws.Cells.Item(12, 10).Value = 3514.1428
$ws.Cells.Item(12, 12).Value = 10542.4284
$ws.Cells.Item(12, 14).Value = -10888.4284
$ws.Cells.Item(54, 8).Value = 1615.3846
$ws.Cells.Item(54, 9).Value = 500
$ws.Cells.Item(54, 10).Value = 3400
$ws.Cells.Item(54, 11).Value = 1500
$ws.Cells.Item(54, 12).Value = 10200
$ws.Cells.Item(54, 13).Value = -941
$ws.Cells.Item(54, 14).Value = -11318
$ws.Cells.Item(131, 8).Value = 17551998
$ws.Cells.Item(131, 10).Value = 11424.615
$ws.Cells.Item(131, 12).Value = 34273.845
$ws.Cells.Item(131, 14).Value = -44353.845
$ws.Cells.Item(132, 8).Value = 1500.6
$ws.Cells.Item(132, 9).Value = 1200.8572
$ws.Cells.Item(132, 10).Value = 2200
$ws.Cells.Item(132, 11).Value = 10807.7148
$ws.Cells.Item(132, 12).Value = 19800
$ws.Cells.Item(132, 13).Value = -8277.7148
$ws.Cells.Item(132, 14).Value = -24860
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2119.0715
$ws.Cells.Item(102, 9).Value = 2189.923
$ws.Cells.Item(102, 11).Value = 2189.923
$ws.Cells.Item(102, 13).Value = -567.9229999999998
$ws.Cells.Item(126, 8).Value = 4120
$ws.Cells.Item(126, 9).Value = 3256
$ws.Cells.Item(126, 10).Value = 4811.2
$ws.Cells.Item(126, 11).Value = 9768
$ws.Cells.Item(126, 12).Value = 14433.6
$ws.Cells.Item(126, 13).Value = -7298
$ws.Cells.Item(126, 14).Value = -19373.6
$ws.Cells.Item(132, 8).Value = 10587
$ws.Cells.Item(132, 9).Value = 9050.777
$ws.Cells.Item(132, 11).Value = 27152.331
$ws.Cells.Item(132, 13).Value = -24622.331
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3311.6487
$ws.Cells.Item(22, 9).Value = 2074.4736
$ws.Cells.Item(22, 10).Value = 4617.5557
$ws.Cells.Item(22, 11).Value = 2074.4736
$ws.Cells.Item(22, 12).Value = 4617.5557
$ws.Cells.Item(22, 13).Value = -1779.4736
$ws.Cells.Item(22, 14).Value = -5207.5557
$ws.Cells.Item(27, 8).Value = 3311.6487
$ws.Cells.Item(27, 9).Value = 2074.4736
$ws.Cells.Item(27, 10).Value = 4617.5557
$ws.Cells.Item(27, 11).Value = 2074.4736
$ws.Cells.Item(27, 12).Value = 4617.5557
$ws.Cells.Item(27, 13).Value = -1967.4736
$ws.Cells.Item(27, 14).Value = -4831.5557
$ws.Cells.Item(46, 8).Value = 4878.8047
$ws.Cells.Item(46, 9).Value = 1349.5
$ws.Cells.Item(46, 11).Value = 1349.5
$ws.Cells.Item(46, 13).Value = -1161.5
$ws.Cells.Item(55, 8).Value = 537.5
$ws.Cells.Item(55, 9).Value = 610.75
$ws.Cells.Item(55, 10).Value = 391
$ws.Cells.Item(55, 11).Value = 610.75
$ws.Cells.Item(55, 12).Value = 391
$ws.Cells.Item(55, 13).Value = -437.75
$ws.Cells.Item(55, 14).Value = -737
$ws.Cells.Item(109, 8).Value = 73820.75
$ws.Cells.Item(109, 10).Value = 73820.75
$ws.Cells.Item(109, 12).Value = 73820.75
$ws.Cells.Item(109, 14).Value = -76594.75
$ws.Cells.Item(122, 8).Value = 4817.88
$ws.Cells.Item(122, 9).Value = 4337
$ws.Cells.Item(122, 11).Value = 13011
$ws.Cells.Item(122, 13).Value = -10561
$ws.Cells.Item(132, 8).Value = 9762.25
$ws.Cells.Item(132, 10).Value = 16524.5
$ws.Cells.Item(132, 12).Value = 49573.5
$ws.Cells.Item(132, 14).Value = -54633.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(24, 8).Value = 52500
$ws.Cells.Item(24, 9).Value = 30000
$ws.Cells.Item(24, 11).Value = 30000
$ws.Cells.Item(24, 13).Value = -29770
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).Value = $null
$ws.Cells.Item(122, 8).Value = 7693.2383
$ws.Cells.Item(122, 9).Value = 7753.222
$ws.Cells.Item(122, 11).Value = 23259.666
$ws.Cells.Item(122, 13).Value = -20809.666
$ws.Cells.Item(132, 8).Value = 6801.7393
$ws.Cells.Item(132, 9).Value = 3826.2222
$ws.Cells.Item(132, 10).Value = 8714.571
$ws.Cells.Item(132, 11).Value = 11478.6666
$ws.Cells.Item(132, 12).Value = 26143.713
$ws.Cells.Item(132, 13).Value = -8948.6666
$ws.Cells.Item(132, 14).Value = -31203.713
